# Commit: "app dashboard updates to add farmers screen and more + backend
# additions" — the documentation-relevant portion of this change appends
# three new REST endpoint rows (getFarmTypes, getCrops, getLocations) to the
# bottom of the single "REST API Endpoints" table.
#
# (The rest of the underlying XML diff is Word's own incidental run
# splitting/merging and proofErr spell-check marker churn around text that
# is otherwise byte-for-byte identical before/after, so there is no visible
# content to reproduce there.)

$d = $word.ActiveDocument
$tbl = $d.Tables.Item(1)

$boilerplate401 = "access token disabled via signout (401) / access token expired (401) / not authorized to access this (401) / invalid token (401)"

$newRows = @(
    @{
        Endpoint = "getFarmTypes";
        Type = "POST";
        Roles = "admin, clerk";
        Headers = "Access-Token";
        Input = "";
        Return = "Farm Types Object List (200) / " + $boilerplate401;
    },
    @{
        Endpoint = "getCrops";
        Type = "POST";
        Roles = "admin, clerk";
        Headers = "Access-Token";
        Input = "";
        Return = "Crop Object List / " + $boilerplate401;
    },
    @{
        Endpoint = "getLocations";
        Type = "POST";
        Roles = "admin, clerk";
        Headers = "Access-Token";
        Input = "";
        Return = "Location Object List / " + $boilerplate401;
    }
)

foreach ($rowData in $newRows) {
    $newRow = $tbl.Rows.Add()
    $newRow.Cells.Item(1).Range.Text = $rowData.Endpoint
    $newRow.Cells.Item(2).Range.Text = $rowData.Type
    $newRow.Cells.Item(3).Range.Text = $rowData.Roles
    $newRow.Cells.Item(4).Range.Text = $rowData.Headers
    if ($rowData.Input -ne "") {
        $newRow.Cells.Item(5).Range.Text = $rowData.Input
    }
    $newRow.Cells.Item(6).Range.Text = $rowData.Return
}

Write-Output ("Table now has " + $tbl.Rows.Count + " rows")
